$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D2:E51 stays text-typed so numeric-looking strings are not
# auto-converted to numbers when assigned via .Value, then restore the
# default (Normal) style so no stray formatting is introduced.
$changedRange = $ws.Range("D2:E51")
$changedRange.NumberFormat = "@"

$ws.Range("D2").Value = "27.229.47"
$ws.Range("E2").Value = "  +0.99%  "
$ws.Range("D3").Value = "1.690.19"
$ws.Range("E3").Value = "  +0.72%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "216.36"
$ws.Range("E5").Value = "  +0.66%  "
$ws.Range("E6").Value = "  +0.84%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "23.18"
$ws.Range("E8").Value = "  +13.72%  "
$ws.Range("E9").Value = "  +4.63%  "
$ws.Range("E10").Value = "  +1.57%  "
$ws.Range("D11").Value = "0.0891"
$ws.Range("E11").Value = "  +0.33%  "
$ws.Range("D12").Value = "1.928.51"
$ws.Range("E12").Value = "  +0.69%  "
$ws.Range("D13").Value = "1.690.72"
$ws.Range("E13").Value = "  -0.12%  "
$ws.Range("E14").Value = "  +2.73%  "
$ws.Range("D15").Value = "0.556"
$ws.Range("E15").Value = "  +5.39%  "
$ws.Range("D16").Value = "67.72"
$ws.Range("E16").Value = "  +3.15%  "
$ws.Range("D17").Value = "27.225.10"
$ws.Range("E17").Value = "  +0.86%  "
$ws.Range("D18").Value = "237.52"
$ws.Range("E18").Value = "  +0.73%  "
$ws.Range("D19").Value = "8.13"
$ws.Range("E19").Value = "  -0.59%  "
$ws.Range("D20").Value = "0.0₃0746"
$ws.Range("E20").Value = "  +1.69%  "
$ws.Range("E21").Value = "  +0.03%  "
$ws.Range("E22").Value = "  +2.86%  "
$ws.Range("D23").Value = "9.64"
$ws.Range("E23").Value = "  +5.18%  "
$ws.Range("E24").Value = "  -1.49%  "
$ws.Range("D25").Value = "147.44"
$ws.Range("E25").Value = "  +0.47%  "
$ws.Range("E26").Value = "  +1.44%  "
$ws.Range("D27").Value = "16.47"
$ws.Range("E27").Value = "  +2.73%  "
$ws.Range("E28").Value = "  +1.17%  "
$ws.Range("E29").Value = "  +0.18%  "
$ws.Range("D30").Value = "0.0504"
$ws.Range("E30").Value = "  +0.80%  "
$ws.Range("E31").Value = "  +1.05%  "
$ws.Range("D32").Value = "3.40"
$ws.Range("E32").Value = "  +2.23%  "
$ws.Range("D33").Value = "1.545.99"
$ws.Range("E33").Value = "  +4.39%  "
$ws.Range("E34").Value = "  +2.19%  "
$ws.Range("E35").Value = "  -0.25%  "
$ws.Range("D36").Value = "0.951"
$ws.Range("E36").Value = "  +4.87%  "
$ws.Range("E37").Value = "  +3.42%  "
$ws.Range("E38").Value = "  -0.59%  "
$ws.Range("E39").Value = "  -0.79%  "
$ws.Range("E40").Value = "  +3.98%  "
$ws.Range("D41").Value = "69.33"
$ws.Range("E41").Value = "  +3.04%  "
$ws.Range("D42").Value = "5.75"
$ws.Range("E42").Value = "  -1.07%  "
$ws.Range("D44").Value = "2.27"
$ws.Range("E44").Value = "  -2.10%  "
$ws.Range("D45").Value = "1.835.71"
$ws.Range("E45").Value = "  +0.86%  "
$ws.Range("E46").Value = "  +0.70%  "
$ws.Range("D47").Value = "91.30"
$ws.Range("E47").Value = "  +1.02%  "
$ws.Range("E48").Value = "  +3.47%  "
$ws.Range("E49").Value = "  +5.44%  "
$ws.Range("D50").Value = "8.28"
$ws.Range("E50").Value = "  +7.39%  "
$ws.Range("E51").Value = "  +2.05%  "

# Restore default styling (remove the quote-prefix/text-format styling marker)
$changedRange.Style = "Normal"
